# Weekly update: insert two new price records for "Poroto verde" (Magnum /
# Sin especificar) at Mercado Mayorista Lo Valledor de Santiago, pushing the
# existing history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 1053 - this shifts the
# old rows 1053:1088 down to 1055:1090 (and carries the D-column date style
# down with them, same as Excel's native Insert behaviour).
$ws.Rows.Item(1053).Insert()
$ws.Rows.Item(1053).Insert()

# --- New row 1053 ---------------------------------------------------------
$ws.Cells.Item(1053, 1).Value  = 6
$ws.Cells.Item(1053, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1053, 3).Value  = "Metropolitana"
$ws.Cells.Item(1053, 4).Value  = 44939
$ws.Cells.Item(1053, 5).Value  = 13
$ws.Cells.Item(1053, 6).Value  = 100112031
$ws.Cells.Item(1053, 7).Value  = "Poroto verde"
$ws.Cells.Item(1053, 8).Value  = "Magnum"
$ws.Cells.Item(1053, 9).Value  = "Primera"
$ws.Cells.Item(1053, 10).Value = 1430
$ws.Cells.Item(1053, 11).Value = 22000
$ws.Cells.Item(1053, 12).Value = 25000
$ws.Cells.Item(1053, 13).Value = 23154
$ws.Cells.Item(1053, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(1053, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1053, 16).Value = 926
$ws.Cells.Item(1053, 17).Value = 25
$ws.Cells.Item(1053, 18).Value = "Hortaliza"

# --- New row 1054 ---------------------------------------------------------
$ws.Cells.Item(1054, 1).Value  = 6
$ws.Cells.Item(1054, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1054, 3).Value  = "Metropolitana"
$ws.Cells.Item(1054, 4).Value  = 44939
$ws.Cells.Item(1054, 5).Value  = 13
$ws.Cells.Item(1054, 6).Value  = 100112031
$ws.Cells.Item(1054, 7).Value  = "Poroto verde"
$ws.Cells.Item(1054, 8).Value  = "Sin especificar"
$ws.Cells.Item(1054, 9).Value  = "Primera"
$ws.Cells.Item(1054, 10).Value = 450
$ws.Cells.Item(1054, 11).Value = 33000
$ws.Cells.Item(1054, 12).Value = 35000
$ws.Cells.Item(1054, 13).Value = 33978
$ws.Cells.Item(1054, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(1054, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1054, 16).Value = 1359
$ws.Cells.Item(1054, 17).Value = 25
$ws.Cells.Item(1054, 18).Value = "Hortaliza"
